# Update gh-pages output data (generated at 456a3b4)
# Bumps "want-to-go" counts (column F) and refreshes one event's cover image
# (column I) across the relevant worksheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (exhibitions) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 70
$ws1.Range("F3").Value = 128
$ws1.Range("F4").Value = 2066
$ws1.Range("F5").Value = 354
$ws1.Range("F6").Value = 601
$ws1.Range("F8").Value = 2066
$ws1.Range("F9").Value = 10604
$ws1.Range("F10").Value = 180
$ws1.Range("F12").Value = 281
$ws1.Range("F14").Value = 416
$ws1.Range("F15").Value = 7488
$ws1.Range("F18").Value = 237
$ws1.Range("F20").Value = 3321

# ---- Sheet "演出" (performances) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("I3").Value = "//i1.hdslb.com/bfs/openplatform/202405/MTs1Gl1Z1715588874037.jpeg"

# ---- Sheet "全部类型" (all types, union of the others) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 70
$ws4.Range("F3").Value = 128
$ws4.Range("F4").Value = 2066
$ws4.Range("F5").Value = 354
$ws4.Range("F6").Value = 601
$ws4.Range("F9").Value = 2066
$ws4.Range("F12").Value = 10604
$ws4.Range("F13").Value = 180
$ws4.Range("F15").Value = 281
$ws4.Range("F17").Value = 416
$ws4.Range("F18").Value = 7488
$ws4.Range("F21").Value = 237
$ws4.Range("F23").Value = 3321
$ws4.Range("I10").Value = "//i1.hdslb.com/bfs/openplatform/202405/MTs1Gl1Z1715588874037.jpeg"
